$p = $ppt.ActivePresentation

# --- 1. Table style change on Slide 5 (table is Shape 2 on that slide) ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{CDE6ACE5-4ACA-430B-A17F-866861F3329C}")

# --- 2. Swap the deck's theme colours ("Integral"/Red Violet -> Office) ---
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

function HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToVbaRgb($officeColors[$i - 1])
}
